$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 7.912706907596791
$ws.Range("D2").Value = 4.156067313357373
$ws.Range("E2").Value = 10.37794237335392
$ws.Range("F2").Value = 63.30968393511017
$ws.Range("G2").Value = 3.788138677660062
$ws.Range("J2").Value = 10.50720820900608
$ws.Range("K2").Value = 18.30029352255165
$ws.Range("M2").Value = 19.05373918289274
$ws.Range("B3").Value = 7.848610166299977
$ws.Range("D3").Value = 4.106801428835364
$ws.Range("E3").Value = 10.41586474912769
$ws.Range("F3").Value = 62.45627643225887
$ws.Range("G3").Value = 3.792548819160817
$ws.Range("J3").Value = 10.5069123178699
$ws.Range("K3").Value = 18.27358556808354
$ws.Range("M3").Value = 19.10555646131287
$ws.Range("B4").Value = 7.810942253039857
$ws.Range("D4").Value = 4.077259667920331
$ws.Range("E4").Value = 10.44074945874019
$ws.Range("F4").Value = 61.93269357720263
$ws.Range("G4").Value = 3.795393414181876
$ws.Range("J4").Value = 10.50756568988588
$ws.Range("K4").Value = 18.26483862089576
$ws.Range("M4").Value = 19.14315016938513
$ws.Range("B5").Value = 7.796032087916601
$ws.Range("D5").Value = 4.065408831346025
$ws.Range("E5").Value = 10.4512933621178
$ws.Range("F5").Value = 61.71960421174053
$ws.Range("G5").Value = 3.796587142362208
$ws.Range("J5").Value = 10.50804189237857
$ws.Range("K5").Value = 18.26319860326296
$ws.Range("M5").Value = 19.15991805460574
$ws.Range("B6").Value = 7.793583264651499
$ws.Range("D6").Value = 4.063452617842469
$ws.Range("E6").Value = 10.4530685439341
$ws.Range("F6").Value = 61.68424246551374
$ws.Range("G6").Value = 3.796787449971718
$ws.Range("J6").Value = 10.5081336442386
$ws.Range("K6").Value = 18.26304247293464
$ws.Range("M6").Value = 19.16278966832975
$ws.Range("B7").Value = 7.810739368773925
$ws.Range("D7").Value = 4.077099071036441
$ws.Range("E7").Value = 10.44089002399239
$ws.Range("F7").Value = 61.92981844154487
$ws.Range("G7").Value = 3.795409373205725
$ws.Range("J7").Value = 10.50757126210182
$ws.Range("K7").Value = 18.26480871261396
$ws.Range("M7").Value = 19.14337045010299
$ws.Range("B8").Value = 7.890265552504042
$ws.Range("D8").Value = 4.138939530187269
$ws.Range("E8").Value = 10.39068635972137
$ws.Range("F8").Value = 63.01544621014634
$ws.Range("G8").Value = 3.789630999831816
$ws.Range("J8").Value = 10.50693285794694
$ws.Range("K8").Value = 18.28949739150378
$ws.Range("M8").Value = 19.07040390755884
$ws.Range("B9").Value = 8.058893475823917
$ws.Range("D9").Value = 4.265375635853896
$ws.Range("E9").Value = 10.30489874884906
$ws.Range("F9").Value = 65.14011811406412
$ws.Range("G9").Value = 3.779378028292558
$ws.Range("J9").Value = 10.5123083766045
$ws.Range("K9").Value = 18.39850025359344
$ws.Range("M9").Value = 18.97335434699354
$ws.Range("B10").Value = 8.189511402480823
$ws.Range("D10").Value = 4.36079506160779
$ws.Range("E10").Value = 10.24954117603298
$ws.Range("F10").Value = 66.68854371857732
$ws.Range("G10").Value = 3.772493270515148
$ws.Range("J10").Value = 10.52030070115537
$ws.Range("K10").Value = 18.51519076550067
$ws.Range("M10").Value = 18.93035896745345
$ws.Range("B11").Value = 8.250170321426717
$ws.Range("D11").Value = 4.40461535991216
$ws.Range("E11").Value = 10.22601311752958
$ws.Range("F11").Value = 67.3881566897322
$ws.Range("G11").Value = 3.769499945592887
$ws.Range("J11").Value = 10.52481452554483
$ws.Range("K11").Value = 18.5761036858401
$ws.Range("M11").Value = 18.91698812241594
$ws.Range("B12").Value = 8.27329940883282
$ws.Range("D12").Value = 4.421256762932778
$ws.Range("E12").Value = 10.21734079457998
$ws.Range("F12").Value = 67.65222575448183
$ws.Range("G12").Value = 3.76838622499521
$ws.Range("J12").Value = 10.52664995497098
$ws.Range("K12").Value = 18.60028183592968
$ws.Range("M12").Value = 18.91281738010174
$ws.Range("B13").Value = 8.268311391636885
$ws.Range("D13").Value = 4.417670810131937
$ws.Range("E13").Value = 10.21919799338783
$ws.Range("F13").Value = 67.59539466856683
$ws.Range("G13").Value = 3.76862520689604
$ws.Range("J13").Value = 10.52624905456319
$ws.Range("K13").Value = 18.59502544309159
$ws.Range("M13").Value = 18.91367589243644
$ws.Range("B14").Value = 8.252070093363722
$ws.Range("D14").Value = 4.405983571834851
$ws.Range("E14").Value = 10.22529488932564
$ws.Range("F14").Value = 67.40990000247945
$ws.Range("G14").Value = 3.769407923388686
$ws.Range("J14").Value = 10.52496300019016
$ws.Range("K14").Value = 18.57807064494352
$ws.Range("M14").Value = 18.91662709227439
$ws.Range("B15").Value = 8.242141938666412
$ws.Range("D15").Value = 4.398830649361954
$ws.Range("E15").Value = 10.22906029082294
$ws.Range("F15").Value = 67.29616222626618
$ws.Range("G15").Value = 3.76988993242925
$ws.Range("J15").Value = 10.52419167788335
$ws.Range("K15").Value = 18.56782966462488
$ws.Range("M15").Value = 18.9185510887778
$ws.Range("B16").Value = 8.185570494940812
$ws.Range("D16").Value = 4.357938609240524
$ws.Range("E16").Value = 10.25111202400448
$ws.Range("F16").Value = 66.64271143605966
$ws.Range("G16").Value = 3.772691667907425
$ws.Range("J16").Value = 10.5200233714548
$ws.Range("K16").Value = 18.51136636003143
$ws.Range("M16").Value = 18.93135754273946
$ws.Range("B17").Value = 8.151169593780688
$ws.Range("D17").Value = 4.332950141765773
$ws.Range("E17").Value = 10.26506330382212
$ws.Range("F17").Value = 66.24049355830572
$ws.Range("G17").Value = 3.774445836036572
$ws.Range("J17").Value = 10.51769106982587
$ws.Range("K17").Value = 18.47872400664203
$ws.Range("M17").Value = 18.94080084952994
$ws.Range("B18").Value = 8.131501015462081
$ws.Range("D18").Value = 4.318617070363544
$ws.Range("E18").Value = 10.27324346741321
$ws.Range("F18").Value = 66.0087092264158
$ws.Range("G18").Value = 3.775467839988274
$ws.Range("J18").Value = 10.51643225000296
$ws.Range("K18").Value = 18.46068696569807
$ws.Range("M18").Value = 18.94681470302531
$ws.Range("B19").Value = 8.124862414103923
$ws.Range("D19").Value = 4.313771307352781
$ws.Range("E19").Value = 10.27603990172558
$ws.Range("F19").Value = 65.93016093165764
$ws.Range("G19").Value = 3.775816119127319
$ws.Range("J19").Value = 10.51602023509162
$ws.Range("K19").Value = 18.45470707581644
$ws.Range("M19").Value = 18.94895081292132
$ws.Range("B20").Value = 8.154819565648655
$ws.Range("D20").Value = 4.33560619075673
$ws.Range("E20").Value = 10.2635620508622
$ws.Range("F20").Value = 66.2833570645479
$ws.Range("G20").Value = 3.774257751954233
$ws.Range("J20").Value = 10.51793079327372
$ws.Range("K20").Value = 18.4821225584186
$ws.Range("M20").Value = 18.93973530511218
$ws.Range("B21").Value = 8.256836400493425
$ws.Range("D21").Value = 4.409415199319132
$ws.Range("E21").Value = 10.22349764967463
$ws.Range("F21").Value = 67.46440895342587
$ws.Range("G21").Value = 3.769177484709225
$ws.Range("J21").Value = 10.52533732324187
$ws.Range("K21").Value = 18.58302063178741
$ws.Range("M21").Value = 18.91573601217802
$ws.Range("B22").Value = 8.324427213704496
$ws.Range("D22").Value = 4.457925683996978
$ws.Range("E22").Value = 10.19869578108527
$ws.Range("F22").Value = 68.231210037838
$ws.Range("G22").Value = 3.765972501988667
$ws.Range("J22").Value = 10.53091308856717
$ws.Range("K22").Value = 18.6554353258629
$ws.Range("M22").Value = 18.90525416936126
$ws.Range("B23").Value = 8.288275290045299
$ws.Range("D23").Value = 4.43201373945981
$ws.Range("E23").Value = 10.21180671652012
$ws.Range("F23").Value = 67.82247491417607
$ws.Range("G23").Value = 3.767672561621485
$ws.Range("J23").Value = 10.52786997984128
$ws.Range("K23").Value = 18.61619931828573
$ws.Range("M23").Value = 18.91037169316056
$ws.Range("B24").Value = 8.153169072513583
$ws.Range("D24").Value = 4.334405286911685
$ws.Range("E24").Value = 10.26424027082062
$ws.Range("F24").Value = 66.2639801543997
$ws.Range("G24").Value = 3.774342742678252
$ws.Range("J24").Value = 10.51782215872158
$ws.Range("K24").Value = 18.48058379993731
$ws.Range("M24").Value = 18.94021521652311
$ws.Range("B25").Value = 8.012019205131603
$ws.Range("D25").Value = 4.230681417234951
$ws.Range("E25").Value = 10.32675609275132
$ws.Range("F25").Value = 64.56682957651904
$ws.Range("G25").Value = 3.782037242918423
$ws.Range("J25").Value = 10.51014423057723
$ws.Range("K25").Value = 18.36255009905466
$ws.Range("M25").Value = 18.99465348322428
